$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header cell format (H1) onto the new header cells so they pick up
# the same shared style (bold font, border, centered/top aligned) used by the
# rest of row 1, instead of creating a brand-new style entry.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data values for the new I (I0) / J (IF) columns, rows 2-23
$values = @(
    8,
    8,
    6,
    10,
    9,
    9,
    9,
    9,
    8,
    8,
    8,
    8,
    8,
    8,
    8,
    8,
    8,
    8,
    8,
    9,
    8,
    7
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i]
    $ws.Cells.Item($row, 10).Value = $values[$i]
}
